$d = $word.ActiveDocument

# Locate the paragraph immediately following the "NEH" heading paragraph.
# It is an otherwise-empty paragraph that contains only the italicized
# full book title "Néhémie". The edit removes this whole paragraph
# (including its paragraph mark), so the "NEH" heading paragraph is
# directly followed by the paragraph that used to come after it.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $r = $p.Range
    if ($r.Italic -and ($r.Text.TrimEnd([char]13, [char]7) -eq "Néhémie")) {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $target.Range.Delete()
}
